# Add a new "parameters" worksheet with parameter/value pairs,
# move it to the end of the workbook (after "units"), and make it
# the active/selected tab.

$wb = $excel.ActiveWorkbook

# Add the new sheet, then move it to the end of the workbook.
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "parameters"
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Move($null, $lastSheet)

# Re-fetch the sheet by name: the reference held before the Move can go stale.
$ws = $wb.Worksheets.Item("parameters")

# Header row.
$ws.Range("A1").Value = "parameter"
$ws.Range("B1").Value = "value"
$ws.Range("A1:B1").Font.Bold = $true
$ws.Range("A1:B1").Font.Italic = $true

# Data rows.
$data = @(
    @("melee_distance", 4.5),
    @("melee_height_difference_threshold", 2),
    @("archer_distance", 4.5),
    @("archer_distance_height_gain", 0.5),
    @("siege_distance", 11),
    @("siege_distance_height_gain", 0.5),
    @("flier_distance", 10),
    @("flier_distance_height_gain", 0.5)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $row++
}

# Column widths to match source formatting.
$ws.Columns.Item(1).ColumnWidth = 34.140625
$ws.Columns.Item(2).ColumnWidth = 5.85546875

# Make the new sheet the active/selected tab.
$ws.Select()
$ws.Range("D31").Select()
